# "Hong fix doi ten" - drop the stray "F"/"G" columns (I & J) that were
# left over from a rename, keep just A-H, and rename the former "G"
# header (now H1) to "Z".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the cells that belonged to the removed I/J columns (this also
# drops the now-unused "F"/"G" shared strings and shrinks the sheet's
# used range back down to A1:H8, matching the target file).
$ws.Range("I1").ClearContents() | Out-Null
$ws.Range("J1").ClearContents() | Out-Null
$ws.Range("H2").ClearContents() | Out-Null
$ws.Range("J2").ClearContents() | Out-Null
$ws.Range("I4").ClearContents() | Out-Null
$ws.Range("J5").ClearContents() | Out-Null
$ws.Range("J7").ClearContents() | Out-Null

# Rename the header that used to read "G" to "Z".
$ws.Range("H1").Value = "Z"

# Leave the selection where the author's saved file shows it.
$ws.Range("I7").Select() | Out-Null
